# Auto-generated edit script applying the data refresh for Linea 141 - LP1912 (12/01/2026)
# Scrape re-run at 11:23:54: updates timestamps/counts, re-sorts/updates rows near
# the scrape boundary, and appends newly observed arrivals.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Sheet "LP1912" ---
$ws1.Range("A2").Value = 'Última actualización: 11:23:54'
$ws1.Range("A3").Value = 'Total filas: 216'
$ws1.Range("A45").Value = '05:20:00'
$ws1.Range("C45").Value = '11_ETCHEVERRY'
$ws1.Range("D45").Value = 116
$ws1.Range("A46").Value = '06:52:23'
$ws1.Range("C46").Value = '16_SANTA ANA'
$ws1.Range("D46").Value = 24
$ws1.Range("A124").Value = '08:39:56'
$ws1.Range("C124").Value = '215C_EL PATO'
$ws1.Range("D124").Value = 62
$ws1.Range("A125").Value = '09:38:04'
$ws1.Range("C125").Value = '16_SANTA ANA'
$ws1.Range("D125").Value = 3
$ws1.Range("C126").Value = '14_ABASTO'
$ws1.Range("A152").Value = '10:57:58'
$ws1.Range("C152").Value = '17_ROMERO'
$ws1.Range("D152").Value = 0
$ws1.Range("A154").Value = '10:28:12'
$ws1.Range("C154").Value = '23_HERNANDEZ'
$ws1.Range("D154").Value = 29
$ws1.Range("A171").Value = '11:23:54'
$ws1.Range("B171").Value = '11:23'
$ws1.Range("C171").Value = '16_SANTA ANA'
$ws1.Range("D171").Value = 0
$ws1.Range("A172").Value = '11:23:54'
$ws1.Range("B172").Value = '11:23'
$ws1.Range("C172").Value = '17_ROMERO'
$ws1.Range("D172").Value = 0
$ws1.Range("A173").Value = '11:23:54'
$ws1.Range("B173").Value = '11:24'
$ws1.Range("C173").Value = '15_ABASTO'
$ws1.Range("D173").Value = 1
$ws1.Range("A174").Value = '11:23:54'
$ws1.Range("B174").Value = '11:25'
$ws1.Range("D174").Value = 2
$ws1.Range("A175").Value = '09:38:04'
$ws1.Range("B175").Value = '11:25'
$ws1.Range("C175").Value = '16_P MOR-SANTA ANA'
$ws1.Range("D175").Value = 107
$ws1.Range("A176").Value = '10:57:58'
$ws1.Range("B176").Value = '11:26'
$ws1.Range("C176").Value = '225_C ROCA-H SUR'
$ws1.Range("D176").Value = 29
$ws1.Range("A177").Value = '10:57:58'
$ws1.Range("B177").Value = '11:26'
$ws1.Range("C177").Value = '23_HERNANDEZ'
$ws1.Range("D177").Value = 29
$ws1.Range("A178").Value = '09:38:04'
$ws1.Range("B178").Value = '11:27'
$ws1.Range("C178").Value = '225_C ROCA-H SUR'
$ws1.Range("D178").Value = 109
$ws1.Range("B179").Value = '11:31'
$ws1.Range("C179").Value = '81_EL PELIGRO'
$ws1.Range("D179").Value = 34
$ws1.Range("A180").Value = '09:38:04'
$ws1.Range("B180").Value = '11:32'
$ws1.Range("C180").Value = '81_EL PELIGRO'
$ws1.Range("D180").Value = 114
$ws1.Range("A181").Value = '11:23:54'
$ws1.Range("B181").Value = '11:34'
$ws1.Range("C181").Value = '23_HERNANDEZ'
$ws1.Range("D181").Value = 11
$ws1.Range("A182").Value = '09:38:04'
$ws1.Range("B182").Value = '11:36'
$ws1.Range("D182").Value = 118
$ws1.Range("A183").Value = '10:28:12'
$ws1.Range("B183").Value = '11:40'
$ws1.Range("C183").Value = '11_ETCHEVERRY'
$ws1.Range("D183").Value = 72
$ws1.Range("A184").Value = '10:57:58'
$ws1.Range("B184").Value = '11:41'
$ws1.Range("C184").Value = '17_ROMERO'
$ws1.Range("D184").Value = 44
$ws1.Range("B185").Value = '11:42'
$ws1.Range("C185").Value = '17_ROMERO'
$ws1.Range("D185").Value = 74
$ws1.Range("B186").Value = '11:43'
$ws1.Range("C186").Value = '10_OLMOS'
$ws1.Range("D186").Value = 46
$ws1.Range("B187").Value = '11:48'
$ws1.Range("C187").Value = '11_ETCHEVERRY'
$ws1.Range("D187").Value = 51
$ws1.Range("A188").Value = '10:57:58'
$ws1.Range("B188").Value = '11:50'
$ws1.Range("C188").Value = '215B_EL PATO'
$ws1.Range("D188").Value = 53
$ws1.Range("A189").Value = '11:23:54'
$ws1.Range("B189").Value = '11:51'
$ws1.Range("C189").Value = '23_HERNANDEZ'
$ws1.Range("D189").Value = 28
$ws1.Range("A190").Value = '10:28:12'
$ws1.Range("B190").Value = '11:51'
$ws1.Range("C190").Value = '10_OLMOS'
$ws1.Range("D190").Value = 83
$ws1.Range("B191").Value = '11:51'
$ws1.Range("C191").Value = '15_ABASTO'
$ws1.Range("D191").Value = 54
$ws1.Range("B192").Value = '11:51'
$ws1.Range("C192").Value = '215B_EL PATO'
$ws1.Range("D192").Value = 83
$ws1.Range("A193").Value = '11:23:54'
$ws1.Range("B193").Value = '11:52'
$ws1.Range("C193").Value = '15_ABASTO'
$ws1.Range("D193").Value = 29
$ws1.Range("A194").Value = '11:23:54'
$ws1.Range("B194").Value = '11:53'
$ws1.Range("C194").Value = '11_ETCHEVERRY'
$ws1.Range("D194").Value = 30
$ws1.Range("A195").Value = '10:57:58'
$ws1.Range("B195").Value = '11:58'
$ws1.Range("C195").Value = '225_GOMEZ'
$ws1.Range("D195").Value = 61
$ws1.Range("B196").Value = '11:59'
$ws1.Range("C196").Value = '225_GOMEZ'
$ws1.Range("D196").Value = 91
$ws1.Range("B197").Value = '12:06'
$ws1.Range("C197").Value = '16_P MOR-SANTA ANA'
$ws1.Range("D197").Value = 69
$ws1.Range("A198").Value = '10:28:12'
$ws1.Range("B198").Value = '12:06'
$ws1.Range("C198").Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Range("D198").Value = 98
$ws1.Range("A199").Value = '10:57:58'
$ws1.Range("B199").Value = '12:06'
$ws1.Range("C199").Value = '14_ABASTO'
$ws1.Range("D199").Value = 69
$ws1.Range("B200").Value = '12:07'
$ws1.Range("C200").Value = '16_P MOR-SANTA ANA'
$ws1.Range("D200").Value = 99
$ws1.Range("B201").Value = '12:10'
$ws1.Range("C201").Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Range("D201").Value = 73
$ws1.Range("A202").Value = '11:23:54'
$ws1.Range("B202").Value = '12:11'
$ws1.Range("C202").Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Range("D202").Value = 48
$ws1.Range("A203").Value = '11:23:54'
$ws1.Range("B203").Value = '12:12'
$ws1.Range("C203").Value = '10_OLMOS'
$ws1.Range("D203").Value = 49
$ws1.Range("A204").Value = '10:57:58'
$ws1.Range("B204").Value = '12:14'
$ws1.Range("C204").Value = '10_OLMOS'
$ws1.Range("D204").Value = 77
$ws1.Range("E204").Value = 'LP1912'
$ws1.Range("A205").Value = '10:28:12'
$ws1.Range("B205").Value = '12:14'
$ws1.Range("C205").Value = '17_ROMERO'
$ws1.Range("D205").Value = 106
$ws1.Range("E205").Value = 'LP1912'
$ws1.Range("A206").Value = '10:28:12'
$ws1.Range("B206").Value = '12:18'
$ws1.Range("C206").Value = '14_ABASTO'
$ws1.Range("D206").Value = 110
$ws1.Range("E206").Value = 'LP1912'
$ws1.Range("A207").Value = '10:57:58'
$ws1.Range("B207").Value = '12:20'
$ws1.Range("C207").Value = '215A_EL PATO'
$ws1.Range("D207").Value = 83
$ws1.Range("E207").Value = 'LP1912'
$ws1.Range("A208").Value = '10:57:58'
$ws1.Range("B208").Value = '12:20'
$ws1.Range("C208").Value = '26_HERNANDEZ'
$ws1.Range("D208").Value = 83
$ws1.Range("E208").Value = 'LP1912'
$ws1.Range("A209").Value = '10:28:12'
$ws1.Range("B209").Value = '12:21'
$ws1.Range("C209").Value = '215A_EL PATO'
$ws1.Range("D209").Value = 113
$ws1.Range("E209").Value = 'LP1912'
$ws1.Range("A210").Value = '10:28:12'
$ws1.Range("B210").Value = '12:21'
$ws1.Range("C210").Value = '26_HERNANDEZ'
$ws1.Range("D210").Value = 113
$ws1.Range("E210").Value = 'LP1912'
$ws1.Range("A211").Value = '10:57:58'
$ws1.Range("B211").Value = '12:29'
$ws1.Range("C211").Value = '17_ROMERO'
$ws1.Range("D211").Value = 92
$ws1.Range("E211").Value = 'LP1912'
$ws1.Range("A212").Value = '10:57:58'
$ws1.Range("B212").Value = '12:36'
$ws1.Range("C212").Value = '27_EL RETIRO'
$ws1.Range("D212").Value = 99
$ws1.Range("E212").Value = 'LP1912'
$ws1.Range("A213").Value = '10:57:58'
$ws1.Range("B213").Value = '12:37'
$ws1.Range("C213").Value = '17_179 Y 38'
$ws1.Range("D213").Value = 100
$ws1.Range("E213").Value = 'LP1912'
$ws1.Range("A214").Value = '11:23:54'
$ws1.Range("B214").Value = '12:38'
$ws1.Range("C214").Value = '17_179 Y 38'
$ws1.Range("D214").Value = 75
$ws1.Range("E214").Value = 'LP1912'
$ws1.Range("A215").Value = '11:23:54'
$ws1.Range("B215").Value = '12:40'
$ws1.Range("C215").Value = '10_OLMOS'
$ws1.Range("D215").Value = 77
$ws1.Range("E215").Value = 'LP1912'
$ws1.Range("A216").Value = '11:23:54'
$ws1.Range("B216").Value = '12:48'
$ws1.Range("C216").Value = '11_ETCHEVERRY'
$ws1.Range("D216").Value = 85
$ws1.Range("E216").Value = 'LP1912'
$ws1.Range("A217").Value = '11:23:54'
$ws1.Range("B217").Value = '12:54'
$ws1.Range("C217").Value = '17_ROMERO'
$ws1.Range("D217").Value = 91
$ws1.Range("E217").Value = 'LP1912'
$ws1.Range("A218").Value = '11:23:54'
$ws1.Range("B218").Value = '13:06'
$ws1.Range("C218").Value = '16_P MOR-SANTA ANA'
$ws1.Range("D218").Value = 103
$ws1.Range("E218").Value = 'LP1912'
$ws1.Range("A219").Value = '11:23:54'
$ws1.Range("B219").Value = '13:13'
$ws1.Range("C219").Value = '215D_EL PATO'
$ws1.Range("D219").Value = 110
$ws1.Range("E219").Value = 'LP1912'
$ws1.Range("A220").Value = '11:23:54'
$ws1.Range("B220").Value = '13:19'
$ws1.Range("C220").Value = '10_OLMOS'
$ws1.Range("D220").Value = 116
$ws1.Range("E220").Value = 'LP1912'
$ws1.Range("A221").Value = '11:23:54'
$ws1.Range("B221").Value = '13:20'
$ws1.Range("C221").Value = '26_HERNANDEZ'
$ws1.Range("D221").Value = 117
$ws1.Range("E221").Value = 'LP1912'

# --- Sheet "LP1912-215" ---
$ws2.Range("A2").Value = 'Última actualización: 11:23:54'
$ws2.Range("A3").Value = 'Total filas: 26'
$ws2.Range("A31").Value = '11:23:54'
$ws2.Range("B31").Value = '13:13'
$ws2.Range("C31").Value = '215D_EL PATO'
$ws2.Range("D31").Value = 110
$ws2.Range("E31").Value = 'LP1912'

# --- Sheet "6203-6173" ---
$ws3.Range("A2").Value = 'Última actualización: 11:23:54'

